$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-1h (E) columns for rows with refreshed market data
# Prices that look like plain numbers are entered with a leading apostrophe so Excel
# keeps them as text (matching the source data, which stores these as strings).

$ws.Range("D2").Value = "26.425.27"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.724.52"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'243.06"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4913"
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("D8").Value = "'0.2614"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").Value = "'0.06204"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "1.718.96"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.07006"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").Value = "'15.46"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "'4.569"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "'0.5993"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "'77.22"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'0.9997"
$ws.Range("D17").Value = "26.430.02"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'0.000007167"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "1.938.60"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'4.488"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "'8.584"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Value = "'5.164"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'138.13"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'15.26"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'107.05"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'1.715"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "'0.07971"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'3.674"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "'0.04540"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'0.9989"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'2.601"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'0.9959"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "'0.6260"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.9267"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").Value = "'1.960"
$ws.Range("E39").Value = "  -6.10%  "
$ws.Range("D40").Value = "'2.389"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.9995"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'0.01486"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'99.86"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'5.334"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").Value = "'0.3846"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "'6.762"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").Value = "'0.1168"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'0.05366"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D51").Value = "'1.234"
$ws.Range("E51").Value = "  -1.30%  "

# Rows 49 and 50 swap: the EnergySwap and Elrond entries exchange places, each with
# refreshed price/volume data, while row 51 (NEARProtocol) keeps its own updated data.
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.730"
$ws.Range("E49").Value = "  -0.95%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'30.14"
$ws.Range("E50").Value = "  -1.95%  "
